# UndoRedoActivityDiagram.pptx edit script
# - Updates the cached "datetimeFigureOut" date placeholder text on every
#   slide layout from 6/7/2018 -> 1/4/2019
# - Renames the "address book" concept to "finance tracker" in the two
#   activity-diagram labels on slide 1
#
# Text edits are applied through TextRange.Characters(start, length).Text
# rather than reassigning the whole TextRange.Text, so that runs/paragraphs
# that are not touched keep their original <a:r>/<a:rPr> boundaries.
#
# Shapes whose TextFrame.AutoSize is "resize shape to fit text"
# (ppAutoSizeShapeToFitText / <a:spAutoFit/>) get re-laid-out by PowerPoint
# whenever their text changes, which can shift the cached Height by a
# fraction of a point. Snapshot the EMU geometry before the edit and
# restore it afterward (with a +0.5 EMU nudge to counter float32
# round-tripping through the Height/Width/Top/Left point properties) so
# the rest of <a:xfrm> stays byte-identical. Shapes that do not auto-fit
# (ppAutoSizeNone) are left completely alone so an untouched <p:spPr/>
# does not grow an explicit <a:xfrm>.

function Get-ShapeGeometryEmu($shape) {
    return @{
        Height = [double]$shape.Height * 12700.0
        Width  = [double]$shape.Width  * 12700.0
        Top    = [double]$shape.Top    * 12700.0
        Left   = [double]$shape.Left   * 12700.0
    }
}

function Set-ShapeGeometryEmu($shape, $geom) {
    $shape.Width  = ($geom.Width  + 0.5) / 12700.0
    $shape.Height = ($geom.Height + 0.5) / 12700.0
    $shape.Top    = ($geom.Top    + 0.5) / 12700.0
    $shape.Left   = ($geom.Left   + 0.5) / 12700.0
}

function Test-AutoFitShape($shape) {
    if (-not $shape.HasTextFrame) {
        return $false
    }
    # ppAutoSizeShapeToFitText = 1
    return ($shape.TextFrame.AutoSize -eq 1)
}

function Replace-RunText($textRange, $oldText, $newText) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        return $false
    }
    $chars = $textRange.Characters($idx + 1, $oldText.Length)
    $chars.Text = $newText
    return $true
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide layouts: refresh the cached date placeholder text.
#    Every layout's "Date Placeholder" (PlaceholderFormat.Type = 16,
#    ppPlaceholderDate) currently shows "6/7/2018"; update it to
#    "1/4/2019" wherever found.
# ---------------------------------------------------------------------
$master = $p.SlideMaster
$layouts = $master.CustomLayouts

for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    $shapes = $layout.Shapes
    for ($si = 1; $si -le $shapes.Count; $si++) {
        $shp = $shapes.Item($si)
        if (-not $shp.HasTextFrame) {
            continue
        }
        $isDatePlaceholder = $false
        if ($shp.Type -eq 14) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        }
        if (-not $isDatePlaceholder) {
            continue
        }
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "6/7/2018") {
            $autoFit = Test-AutoFitShape $shp
            if ($autoFit) {
                $geom = Get-ShapeGeometryEmu $shp
            }
            $tr.Text = "1/4/2019"
            if ($autoFit) {
                Set-ShapeGeometryEmu $shp $geom
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1: rename "address book" -> "finance tracker" in the
#    activity-diagram callouts.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes

for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    if (-not $shp.HasTextFrame) {
        continue
    }
    $tr = $shp.TextFrame.TextRange
    $full = $tr.Text

    $autoFit = Test-AutoFitShape $shp

    if ($full -eq "[command commits address book]") {
        if ($autoFit) {
            $geom = Get-ShapeGeometryEmu $shp
        }
        Replace-RunText $tr "command commits address book]" "command commits finance tracker]" | Out-Null
        if ($autoFit) {
            Set-ShapeGeometryEmu $shp $geom
        }
    }
    elseif ($full -eq "Purge redundant states and then save address book to addressBookStateList ") {
        if ($autoFit) {
            $geom = Get-ShapeGeometryEmu $shp
        }
        Replace-RunText $tr "Purge redundant states and then save address book to " "Purge redundant states and then save finance tracker to " | Out-Null
        Replace-RunText $tr "addressBookStateList" "financeTrackerStateList" | Out-Null
        if ($autoFit) {
            Set-ShapeGeometryEmu $shp $geom
        }
    }
}
